$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.442.22"
$ws.Range("E2").Value = "  +4.84%  "

$ws.Range("D3").Value = "1.590.36"
$ws.Range("E3").Value = "  +1.66%  "

$ws.Range("E4").Value = "  -0.20%  "

$ws.Range("D5").Value = "'214.60"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.99%  "

$ws.Range("D6").Value = "'0.497"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.34%  "

$ws.Range("E7").Value = "  -0.10%  "

$ws.Range("D8").Value = "'23.98"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +8.74%  "

$ws.Range("E9").Value = "  +0.94%  "

$ws.Range("E10").Value = "  +0.59%  "

$ws.Range("E11").Value = "  +2.21%  "

$ws.Range("D12").Value = "1.816.20"
$ws.Range("E12").Value = "  +1.33%  "

$ws.Range("D13").Value = "1.605.99"
$ws.Range("E13").Value = "  +3.90%  "

$ws.Range("E14").Value = "  +0.47%  "

$ws.Range("D15").Value = "'0.534"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +3.05%  "

$ws.Range("D16").Value = "28.434.08"
$ws.Range("E16").Value = "  +4.24%  "

$ws.Range("D17").Value = "'63.10"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.77%  "

$ws.Range("D18").Value = "'232.74"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +7.34%  "

$ws.Range("E19").Value = "  +1.01%  "

$ws.Range("E20").Value = "  +0.29%  "

$ws.Range("E21").Value = "  -0.40%  "

$ws.Range("E22").Value = "  -0.39%  "

$ws.Range("D23").Value = "'9.43"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.56%  "

$ws.Range("E24").Value = "  +0.48%  "

$ws.Range("E25").Value = "  -0.58%  "

$ws.Range("D26").Value = "'15.27"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.62%  "

$ws.Range("D27").Value = "'6.61"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.27%  "

$ws.Range("E28").Value = "  +1.16%  "

$ws.Range("E29").Value = "  -0.24%  "

$ws.Range("E30").Value = "  +0.58%  "

$ws.Range("E31").Value = "  +0.85%  "

$ws.Range("E32").Value = "  +0.25%  "

$ws.Range("D33").Value = "'3.17"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.32%  "

$ws.Range("D34").Value = "1.419.83"
$ws.Range("E34").Value = "  -1.45%  "

$ws.Range("E35").Value = "  -0.91%  "

$ws.Range("E36").Value = "  -5.23%  "

$ws.Range("E37").Value = "  -0.35%  "

$ws.Range("E38").Value = "  +0.64%  "

$ws.Range("E39").Value = "  +5.29%  "

$ws.Range("E40").Value = "  +2.14%  "

$ws.Range("E41").Value = "  +1.56%  "

$ws.Range("D42").Value = "'5.77"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.59%  "

$ws.Range("E43").Value = "  -0.40%  "

$ws.Range("E44").Value = "  -2.09%  "

$ws.Range("D45").Value = "'1.84"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +6.26%  "

$ws.Range("D46").Value = "'64.64"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.30%  "

$ws.Range("D47").Value = "1.728.40"
$ws.Range("E47").Value = "  +1.36%  "

$ws.Range("D48").Value = "'87.64"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.99%  "

$ws.Range("E49").Value = "  +11.65%  "

$ws.Range("E50").Value = "  -0.60%  "

$ws.Range("D51").Value = "'39.50"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +16.86%  "
